# Commit: "Assign biomass to 100% guaranteed dispatch"
#
# The BGDPbES sheet lists, in column B, the BAU guaranteed-dispatch
# percentage (0-1) for each electricity source named in column A; columns
# C:AK repeat that value out to 2050 via shared formulas ("=$B<row>").
# Row 9 is "biomass" (A9 -> shared string "biomass"). The edit sets its
# guaranteed-dispatch percentage from 0 to 1 (100%); Excel's recalculation
# then propagates that 1 across the shared formulas in C9:AK9.
#
# The author's selection/active-sheet state when the file was saved also
# shifted (they had just edited B9 on BGDPbES, then the workbook was left
# with "About" as the active/selected tab) - replicate that too.

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("BGDPbES")
$wsData.Activate() | Out-Null

# Core data edit: biomass (row 9) guaranteed dispatch -> 100%.
$wsData.Range("B9").Value = 1

# Leave the selection on the edited cell, matching the saved file.
$wsData.Range("B9").Select() | Out-Null

# The workbook was saved with the "About" sheet as the active tab.
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate() | Out-Null
